$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 with new data
$ws.Range("A2").Value = "olganovikov53583@gmail.com"
$ws.Range("B2").Value = "koh93hwlz3ilw7"
$ws.Range("C2").Value = "7knkeepdontity1999zu@aol.com"

$ws.Range("A3").Value = "rozaskiara882@gmail.com"
$ws.Range("B3").Value = "gwpgp6nd6qmkc"
$ws.Range("C3").Value = "tdcounandarse2011c5@yahoo.ca"

$ws.Range("A4").Value = "thidphymee@gmail.com"
$ws.Range("B4").Value = "n0z3ba7rgn6vfo"
$ws.Range("C4").Value = "pv2roychrisicke@divermail.com"

# Delete rows 5-10 (remove the extra rows)
$ws.Range("A5:C10").EntireRow.Delete()
